# Commit: "added element to be clicked and a click script for javascript page"
# The "Data" sheet holds the list of element identifiers (one per row, column A)
# that the automation script clicks through; this edit appends the new batch of
# 164 identifiers that were generated for the new page/script.

$wb = $excel.ActiveWorkbook

$ws   = $wb.Worksheets.Item("Data")
$ws2  = $wb.Worksheets.Item("Database")

$values = @(
    "+3620-796",
    "+593-1194",
    "+6448-633",
    "+17294-220",
    "+1333-594",
    "+17297-220",
    "+21998-106",
    "+22046-182",
    "+19466-646",
    "+552-1087",
    "+173035-33",
    "+754-1480",
    "+17292-220",
    "+814-1445",
    "+817-1445",
    "+813-1445",
    "+17290-220",
    "+816-1445",
    "+173032-33",
    "+173031-33",
    "+22056-182",
    "+788-1402",
    "+15015-335",
    "+17291-220",
    "+17295-220",
    "+17296-220",
    "+12542-314",
    "+811-1445",
    "+26969-232",
    "+30771-65",
    "+17293-220",
    "+17301-220",
    "+3721-696",
    "+187088-77",
    "+169-1530",
    "+5667-281",
    "+2016-705",
    "+0035-1583",
    "+619-577",
    "+1087-1222",
    "+26974-232",
    "+2445-1369",
    "+26979-232",
    "+401-1454",
    "+26978-232",
    "+19923-115",
    "+26972-232",
    "+242-1524",
    "+162-2197",
    "+1858-1321",
    "+916-1429",
    "+1738-1202",
    "+1512-886",
    "+13083-206",
    "+5758-694",
    "+187051-77",
    "+186995-77",
    "+216-1496",
    "+187052-77",
    "+471-1345",
    "186996-77",
    "+3694-687",
    "+187061-77",
    "+442-598",
    "+187058-77",
    "+187073-77",
    "+6396-667",
    "+187055-77",
    "+12310-645",
    "+6830-370",
    "+3517-1088",
    "+192-2175",
    "+239-1506",
    "+662-1341",
    "+1368-2060",
    "+178-1728",
    "+45533-147",
    "+2258-656",
    "+12632-105",
    "+1117-578",
    "+0042-1106",
    "+526-2016",
    "+565-1711",
    "+3443-340",
    "+3542-543",
    "+1532-708",
    "+179-1545",
    "+2051-780",
    "+187068-77",
    "+2443-741",
    "1197-1237",
    "+162-1529",
    "+750-1480",
    "+13949-138",
    "+28359-53",
    "+19828-459",
    "+18368-112",
    "+920-1429",
    "+983-2098",
    "+402-1454",
    "+2964-1211",
    "+13948-138",
    "+0046-1589",
    "+1390-947",
    "+917-1429",
    "+13947-138",
    "+499-1564",
    "+26968-232",
    "+12436-321",
    "+10666-443",
    "+2496-754",
    "+2993-168",
    "+2930-484",
    "+604-1198",
    "+3423-707",
    "+3614-695",
    "+225-1559",
    "+7023-319",
    "+1606-1258",
    "+8832-233",
    "+919-1429",
    "+2764-686",
    "+195-1523",
    "+10670-443",
    "+1350-1071",
    "+163-1532",
    "+1088-828",
    "+2260-656",
    "+403-1454",
    "+3854-377",
    "+0021-1592",
    "+38195-31",
    "+616-2091",
    "+4037-768",
    "+999-859",
    "+406-1454",
    "+1979-651",
    "+12637-105",
    "+6520-664",
    "+202-868",
    "+4248-414",
    "+880-995",
    "+203-868",
    "+12599-410",
    "+4409-632",
    "+2587-1244",
    "+0086-2179",
    "+6520-417",
    "+6519-417",
    "+399-1454",
    "+290-2144",
    "193-1523",
    "+719-2033",
    "+2155-1317",
    "+0069-2200",
    "+12096-292",
    "+3095-790",
    "+405-1454",
    "+0085-2179",
    "+404-1454",
    "+12095-292",
    "+0049-1707",
    "+3088-784",
    "+607-1198",
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match the saved view state: "Data" becomes the visible/active sheet with the
# last-entered row selected, while "Database" keeps its own prior selection but
# is no longer the tab shown on open.
$ws2.Activate()
$ws2.Range("D2:D307").Select()

$ws.Activate()
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 143
    $win.ScrollColumn = 1
} catch {}
$ws.Range("A165").Select()
